$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.246919631958008
$ws.Range("B1").Value = 4.157437801361084
$ws.Range("C1").Value = 1.989984631538391
$ws.Range("D1").Value = 1.512258648872375
$ws.Range("E1").Value = 1.345625400543213
